$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.875.35"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "1.809.74"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'309.76"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D7").Value = "'0.4646"
$ws.Range("E7").Value = "  +4.02%  "
$ws.Range("D8").Value = "'0.3704"
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").Value = "'0.07357"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "'0.8769"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "1.823.52"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "'5.362"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").Value = "'6.522"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("D15").Value = "'91.75"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "'0.07046"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'0.000008694"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D21").Value = "26.875.45"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").Value = "'5.319"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("D24").Value = "2.012.27"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'18.41"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "'2.159"
$ws.Range("E28").Value = "  -4.17%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "'116.06"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "'0.08909"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "'0.7539"
$ws.Range("E32").Value = "  -5.13%  "
$ws.Range("D33").Value = "'1.158"
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.466"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.922"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "'1.101"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'0.01967"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "'0.05255"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").Value = "'2.422"
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("D41").Value = "'2.929"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("D42").Value = "'0.5322"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "'7.172"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").Value = "'0.1664"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("D45").Value = "'8.497"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").Value = "'0.4984"
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("D47").Value = "'10.29"
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "'103.57"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").Value = "'0.06293"
$ws.Range("E51").Value = "  -1.44%  "
